$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Board")

# Insert a new row 32 (pushing the old border/footer rows down to 33/34),
# carrying the formatting down from the row above it.
$ws.Rows("32:32").Insert()
$ws.Range("A31:E31").Copy()
$ws.Range("A32:E32").PasteSpecial(-4122)

# New Gantt data point
$ws.Range("A32").Value = 45431
$ws.Range("B32").Value = 21
$ws.Range("C32").Value = 45

# Extend the D/E formulas down into the new row.
$ws.Range("D32").Formula = "=B32+C32"
$ws.Range("E32").Formula = "=C32/D32"

Write-Output "done"
